$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.955.40'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').Value = '2.636.58'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Formula = "'596.21"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').Formula = "'155.40"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.52%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('E9').Value = '  +5.71%  '
$ws.Range('D10').Formula = "'0.399"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('D13').Formula = "'29.06"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.36%  '
$ws.Range('D14').Formula = "'0.0000188"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +20.49%  '
$ws.Range('D15').Value = '3.113.80'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = '64.857.44'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Formula = "'12.55"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.52%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.538.64'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').Formula = "'4.78"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').Formula = "'351.92"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Formula = "'7.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.29%  '
$ws.Range('D23').Formula = "'68.01"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Formula = "'9.53"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.37%  '
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').Formula = "'8.11"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').Formula = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '0.0₃0937'
$ws.Range('E30').Value = '  +8.90%  '
$ws.Range('E31').Value = '  +3.23%  '
$ws.Range('D32').Formula = "'509.90"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.90%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  +6.63%  '
$ws.Range('D35').Formula = "'6.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('E36').Value = '  +2.37%  '
$ws.Range('D37').Formula = "'164.17"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Formula = "'20.17"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.29%  '
$ws.Range('D39').Formula = "'2.00"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.80%  '
$ws.Range('D40').Formula = "'1.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Formula = "'42.24"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.01%  '
$ws.Range('D43').Formula = "'165.08"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('E45').Value = '  +2.91%  '
$ws.Range('D46').Formula = "'22.99"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('E48').Value = '  +2.96%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Formula = "'0.0981"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').Formula = "'19.32"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.58%  '
